$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block (rows 5-12) grows from 8 blank rows to 11 blank rows so
# that the footer (previously rows 13-16) ends up at rows 16-19, keeping
# its original formatting intact.
$ws.Rows("13:15").Insert()
$ws.Rows("13:15").RowHeight = 18

# Full, updated requirements list (rows 2-14).
$data = @(
    @("Listar produtos",                                                          "x", "",  "RF001"),
    @("Visualizar detalhes do produto",                                           "x", "",  "RF002"),
    @("Cadastrar dados do aluno",                                                  "x", "",  "RF003"),
    @("Visualizar pedidos",                                                       "",  "x", "RF004"),
    @("Efetuar login",                                                            "",  "x", "RF005"),
    @("Cadastrar produto",                                                        "",  "x", "RF006"),
    @("Editar produto",                                                           "",  "x", "RF007"),
    @("Definir disponível ou indisponível",                                       "",  "x", "RF008"),
    @("Definir poucas unidades",                                                  "",  "x", "RF009"),
    @("Exibir informações de contato",                                            "x", "",  "RNF001"),
    @("Exibir mensagem de confirmação do pedido",                                 "x", "",  "RNF002"),
    @("Somente administrador pode acessar o sistema web",                         "",  "x", "RN001"),
    @("Disponibilizar o aplicativo nas plataformas Windows Phone, Android e iOs", "x", "",  "RNF003")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}

$ws.Range("A13").Select()
